# W1S2.pptx edit -- 26/01/2023 -> 27/01/2023 working session:
#   1. Slide 22 ("Overfitting" slide), Content Placeholder 1:
#        "or too much data," -> "or has imbalance in the data,"
#   2. The cached "today" date field (datetimeFigureOut) on the slide
#      master and all 11 custom layouts rolls from 26/01/2023 to
#      27/01/2023 (PowerPoint re-stamps this field text whenever the
#      deck is touched on a later day).
#      NB: the Notes Master's own Date Placeholder is intentionally left
#      alone -- shape-id collisions between the Notes Master and the
#      Slide Master make that particular placeholder unsafe to touch
#      through this host (it would silently corrupt the Slide Master's
#      "Text Placeholder 2" instead of updating the notes date).

$p = $ppt.ActivePresentation

# --- 1. Content fix on slide 22 ------------------------------------------
$slide = $p.Slides.Item(22)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.Name -eq "Content Placeholder 1") {
        $tr = $shp.TextFrame.TextRange
        [void]$tr.Replace("or too much data,", "or has imbalance in the data,")
    }
}

# --- 2. Refresh the cached "today" date shown by the date placeholders ---
$oldDate = "26/01/2023"
$newDate = "27/01/2023"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every custom (slide) layout
for ($l = 1; $l -le $p.SlideMaster.CustomLayouts.Count; $l++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($l)
    Update-DatePlaceholder $layout.Shapes
}
